$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update descriptions in column J for Construcción Iteración 2 (C201-C205)
$ws.Range("J2").Value = "Informe de Revisión tecnica formal."
$ws.Range("J3").Value = "Consultar Escenario"
$ws.Range("J4").Value = "Asignar Escenario"
$ws.Range("J5").Value = "Adjuntar Documentos"
$ws.Range("J6").Value = "Comentar Escenario"
